# Apply "new data (RECOVERY Bari)" edit:
#  - Rename study "RECOVERY" (row 13) to "RECOVERY Toci"
#  - Append two new study rows (27, 28) for "RECOVERY Bari" and
#    "RECOVERY Bari (No Toci)" sourced from "Gray literature"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 27: RECOVERY Bari
$ws.Range("A27").Value = "RECOVERY Bari"
$ws.Range("B27").Value = "Gray literature"
$ws.Range("C27").Value = "baricitinib"
$ws.Range("D27").Value = 487
$ws.Range("E27").Value = 3962
$ws.Range("F27").Value = 523
$ws.Range("G27").Value = 3809

# Rename the existing "RECOVERY" study row to "RECOVERY Toci"
$ws.Range("A13").Value = "RECOVERY Toci"

# New row 28: RECOVERY Bari (No Toci)
$ws.Range("A28").Value = "RECOVERY Bari (No Toci)"
$ws.Range("B28").Value = "Gray literature"
$ws.Range("C28").Value = "baricitinib"
$ws.Range("D28").Value = 306
$ws.Range("E28").Value = 2623
$ws.Range("F28").Value = 310
$ws.Range("G28").Value = 2525

# Column A needs to widen to fit the new, longer study names (bestFit)
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Leave the selection where the editor's cursor ended up
$ws.Range("I27").Select() | Out-Null
